$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Reorder some country name labels (shared-string swap in the source diff) ---
# Santa Lucia <-> Nueva Caledonia
$ws.Range("A198").Value = "Santa Lucia"
$ws.Range("A199").Value = "Nueva Caledonia"

# Groenlandia <-> Seychelles
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Seychelles"

# Sahara Occidental rotates to the front of the Bonaire/San Bartolome group
$ws.Range("A214").Value = "Sahara Occidental"
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "San Bartolome"

# --- Update the "last updated" timestamp label ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 22:05"

# --- Update statistic counts for several countries ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1680400
$ws.Range("C4").Value = 13572
$ws.Range("D4").Value = 451207
$ws.Range("E4").Value = 1130063
$ws.Range("G4").Value = 447
$ws.Range("H4").Value = 99130

# Row 5: Brasil
$ws.Range("B5").Value = 354460
$ws.Range("C5").Value = 7062
$ws.Range("E5").Value = 189461
$ws.Range("G5").Value = 399
$ws.Range("H5").Value = 22412

# Row 11: Alemania
$ws.Range("B11").Value = 180250
$ws.Range("C11").Value = 264
$ws.Range("E11").Value = 11579

# Row 16: Canada
$ws.Range("B16").Value = 84657
$ws.Range("C16").Value = 1036
$ws.Range("D16").Value = 43930
$ws.Range("E16").Value = 34303

# Row 175: Malaui
$ws.Range("B175").Value = 83
$ws.Range("C175").Value = 1
$ws.Range("D175").Value = 33
$ws.Range("E175").Value = 46
